$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.004754781723022
$ws.Range("B1").Value = 2.11422324180603
$ws.Range("C1").Value = 6.672050952911377
$ws.Range("D1").Value = 1.855688333511353
$ws.Range("E1").Value = 1.370176672935486
